# Update Wars Analysis data: append 10 new "EKU" war rows to the WarsAnalysis
# table on the "Wars Analysis" sheet (rows 193-202), mirroring the
# WarNum/BU_Score/Opp_Score/.../HomeGame pattern used by the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wars Analysis")
$lo = $ws.ListObjects.Item("WarsAnalysis")

# New game data for the EKU wars (WarNum, BU_Score, Opp_Score, GameWon,
# WarLabel, Half, GameOrder, HomeGame).
$newGames = @(
    @{ WarNum=1;  BU=7;  Opp=7;  GameWon=0; Label="20-16"; Half=1; Order=20; Home=0 },
    @{ WarNum=2;  BU=2;  Opp=6;  GameWon=0; Label="16-12"; Half=1; Order=20; Home=0 },
    @{ WarNum=3;  BU=5;  Opp=12; GameWon=0; Label="12-8";  Half=1; Order=20; Home=0 },
    @{ WarNum=4;  BU=2;  Opp=8;  GameWon=0; Label="8-4";   Half=1; Order=20; Home=0 },
    @{ WarNum=5;  BU=8;  Opp=7;  GameWon=0; Label="4-0";   Half=1; Order=20; Home=0 },
    @{ WarNum=6;  BU=11; Opp=3;  GameWon=0; Label="20-16"; Half=2; Order=20; Home=0 },
    @{ WarNum=7;  BU=8;  Opp=15; GameWon=0; Label="16-12"; Half=2; Order=20; Home=0 },
    @{ WarNum=8;  BU=10; Opp=12; GameWon=0; Label="12-8";  Half=2; Order=20; Home=0 },
    @{ WarNum=9;  BU=9;  Opp=10; GameWon=0; Label="8-4";   Half=2; Order=20; Home=0 },
    @{ WarNum=10; BU=7;  Opp=9;  GameWon=0; Label="4-0";   Half=2; Order=20; Home=0 }
)

$startRow = $lo.Range.Rows.Count + $lo.Range.Row - 1 + 1  # first blank row after the table header/body

foreach ($game in $newGames) {
    $newRow = $lo.ListRows.Add()
    $r = $lo.Range.Row + $lo.Range.Rows.Count - 1

    $ws.Cells.Item($r, 1).Value = $game.WarNum
    $ws.Cells.Item($r, 2).Value = $game.BU
    $ws.Cells.Item($r, 3).Value = $game.Opp
    $ws.Cells.Item($r, 4).Formula = "=B$r-C$r"
    $ws.Cells.Item($r, 5).Formula = "=IF(WarsAnalysis[[#This Row],[ScoreDiff]]>0,1,0)"
    $ws.Cells.Item($r, 6).Value = $game.GameWon
    $ws.Cells.Item($r, 7).Value = "EKU"
    $ws.Cells.Item($r, 8).Value = 1
    $ws.Cells.Item($r, 9).NumberFormat = "@"
    $ws.Cells.Item($r, 9).Value = $game.Label
    $ws.Cells.Item($r, 10).Value = $game.Half
    $ws.Cells.Item($r, 11).Value = $game.Order
    $ws.Cells.Item($r, 12).Value = $game.Home
}

# Refresh the pivot tables/cache against the (unchanged) WarByLineup source.
$wb.RefreshAll()

# Leave the selection the way the author left it after typing the new rows.
$ws.Range("B203").Select()
